# Update NHL/WHL playoff predictions for rows 379-412 (sheet "Sheet1")
# - Fills in Actual/Correct results for games that have now been played
#   (Fri Mar 21, Sat Mar 22, Sun Mar 23, 2025 games)
# - Adds new rows for the next round of predictions (Fri Mar 28, 2025)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A382").Value2 = 1021928
$ws.Range("B382").Value2 = 'Fri, Mar 21, 2025'
$ws.Range("C382").Value2 = 'Brandon Wheat Kings'
$ws.Range("D382").Value2 = 'Regina Pats'
$ws.Range("E382").Value2 = 'Brandon Wheat Kings'
$ws.Range("F382").Value2 = 'Brandon Wheat Kings'
$ws.Range("G382").Value2 = 1

$ws.Range("A383").Value2 = 1021931
$ws.Range("B383").Value2 = 'Fri, Mar 21, 2025'
$ws.Range("C383").Value2 = 'Lethbridge Hurricanes'
$ws.Range("D383").Value2 = 'Edmonton Oil Kings'
$ws.Range("E383").Value2 = 'Lethbridge Hurricanes'
$ws.Range("F383").Value2 = 'Lethbridge Hurricanes'
$ws.Range("G383").Value2 = 1

$ws.Range("A384").Value2 = 1021932
$ws.Range("B384").Value2 = 'Fri, Mar 21, 2025'
$ws.Range("C384").Value2 = 'Swift Current Broncos'
$ws.Range("D384").Value2 = 'Moose Jaw Warriors'
$ws.Range("E384").Value2 = 'Swift Current Broncos'
$ws.Range("F384").Value2 = 'Swift Current Broncos'
$ws.Range("G384").Value2 = 1

$ws.Range("A385").Value2 = 1021935
$ws.Range("B385").Value2 = 'Fri, Mar 21, 2025'
$ws.Range("C385").Value2 = 'Saskatoon Blades'
$ws.Range("D385").Value2 = 'Prince Albert Raiders'
$ws.Range("E385").Value2 = 'Saskatoon Blades'
$ws.Range("F385").Value2 = 'Prince Albert Raiders'
$ws.Range("G385").Value2 = 0

$ws.Range("A386").Value2 = 1021930
$ws.Range("B386").Value2 = 'Fri, Mar 21, 2025'
$ws.Range("C386").Value2 = 'Kamloops Blazers'
$ws.Range("D386").Value2 = 'Vancouver Giants'
$ws.Range("E386").Value2 = 'Kamloops Blazers'
$ws.Range("F386").Value2 = 'Vancouver Giants'
$ws.Range("G386").Value2 = 0

$ws.Range("A387").Value2 = 1021933
$ws.Range("B387").Value2 = 'Fri, Mar 21, 2025'
$ws.Range("C387").Value2 = 'Portland Winterhawks'
$ws.Range("D387").Value2 = 'Tri-City Americans'
$ws.Range("E387").Value2 = 'Portland Winterhawks'
$ws.Range("F387").Value2 = 'Tri-City Americans'
$ws.Range("G387").Value2 = 0

$ws.Range("A388").Value2 = 1021934
$ws.Range("B388").Value2 = 'Fri, Mar 21, 2025'
$ws.Range("C388").Value2 = 'Prince George Cougars'
$ws.Range("D388").Value2 = 'Victoria Royals'
$ws.Range("E388").Value2 = 'Victoria Royals'
$ws.Range("F388").Value2 = 'Victoria Royals'
$ws.Range("G388").Value2 = 1

$ws.Range("A389").Value2 = 1021929
$ws.Range("B389").Value2 = 'Fri, Mar 21, 2025'
$ws.Range("C389").Value2 = 'Everett Silvertips'
$ws.Range("D389").Value2 = 'Wenatchee Wild'
$ws.Range("E389").Value2 = 'Everett Silvertips'
$ws.Range("F389").Value2 = 'Everett Silvertips'
$ws.Range("G389").Value2 = 1

$ws.Range("A390").Value2 = 1021936
$ws.Range("B390").Value2 = 'Fri, Mar 21, 2025'
$ws.Range("C390").Value2 = 'Spokane Chiefs'
$ws.Range("D390").Value2 = 'Seattle Thunderbirds'
$ws.Range("E390").Value2 = 'Spokane Chiefs'
$ws.Range("F390").Value2 = 'Seattle Thunderbirds'
$ws.Range("G390").Value2 = 0

$ws.Range("A391").Value2 = 1021942
$ws.Range("B391").Value2 = 'Sat, Mar 22, 2025'
$ws.Range("C391").Value2 = 'Regina Pats'
$ws.Range("D391").Value2 = 'Brandon Wheat Kings'
$ws.Range("E391").Value2 = 'Brandon Wheat Kings'
$ws.Range("F391").Value2 = 'Brandon Wheat Kings'
$ws.Range("G391").Value2 = 1

$ws.Range("A392").Value2 = 1021938
$ws.Range("B392").Value2 = 'Sat, Mar 22, 2025'
$ws.Range("C392").Value2 = 'Medicine Hat Tigers'
$ws.Range("D392").Value2 = 'Lethbridge Hurricanes'
$ws.Range("E392").Value2 = 'Medicine Hat Tigers'
$ws.Range("F392").Value2 = 'Medicine Hat Tigers'
$ws.Range("G392").Value2 = 1

$ws.Range("A393").Value2 = 1021939
$ws.Range("B393").Value2 = 'Sat, Mar 22, 2025'
$ws.Range("C393").Value2 = 'Prince Albert Raiders'
$ws.Range("D393").Value2 = 'Saskatoon Blades'
$ws.Range("E393").Value2 = 'Saskatoon Blades'
$ws.Range("F393").Value2 = 'Prince Albert Raiders'
$ws.Range("G393").Value2 = 0

$ws.Range("A394").Value2 = 1021940
$ws.Range("B394").Value2 = 'Sat, Mar 22, 2025'
$ws.Range("C394").Value2 = 'Prince George Cougars'
$ws.Range("D394").Value2 = 'Victoria Royals'
$ws.Range("E394").Value2 = 'Victoria Royals'
$ws.Range("F394").Value2 = 'Prince George Cougars'
$ws.Range("G394").Value2 = 0

$ws.Range("A395").Value2 = 1021941
$ws.Range("B395").Value2 = 'Sat, Mar 22, 2025'
$ws.Range("C395").Value2 = 'Red Deer Rebels'
$ws.Range("D395").Value2 = 'Calgary Hitmen'
$ws.Range("E395").Value2 = 'Calgary Hitmen'
$ws.Range("F395").Value2 = 'Calgary Hitmen'
$ws.Range("G395").Value2 = 1

$ws.Range("A396").Value2 = 1021944
$ws.Range("B396").Value2 = 'Sat, Mar 22, 2025'
$ws.Range("C396").Value2 = 'Moose Jaw Warriors'
$ws.Range("D396").Value2 = 'Swift Current Broncos'
$ws.Range("E396").Value2 = 'Swift Current Broncos'
$ws.Range("F396").Value2 = 'Moose Jaw Warriors'
$ws.Range("G396").Value2 = 0

$ws.Range("A397").Value2 = 1021946
$ws.Range("B397").Value2 = 'Sat, Mar 22, 2025'
$ws.Range("C397").Value2 = 'Wenatchee Wild'
$ws.Range("D397").Value2 = 'Everett Silvertips'
$ws.Range("E397").Value2 = 'Everett Silvertips'
$ws.Range("F397").Value2 = 'Everett Silvertips'
$ws.Range("G397").Value2 = 1

$ws.Range("A398").Value2 = 1021937
$ws.Range("B398").Value2 = 'Sat, Mar 22, 2025'
$ws.Range("C398").Value2 = 'Kelowna Rockets'
$ws.Range("D398").Value2 = 'Kamloops Blazers'
$ws.Range("E398").Value2 = 'Kamloops Blazers'
$ws.Range("F398").Value2 = 'Kelowna Rockets'
$ws.Range("G398").Value2 = 0

$ws.Range("A399").Value2 = 1021943
$ws.Range("B399").Value2 = 'Sat, Mar 22, 2025'
$ws.Range("C399").Value2 = 'Seattle Thunderbirds'
$ws.Range("D399").Value2 = 'Portland Winterhawks'
$ws.Range("E399").Value2 = 'Seattle Thunderbirds'
$ws.Range("F399").Value2 = 'Seattle Thunderbirds'
$ws.Range("G399").Value2 = 1

$ws.Range("A400").Value2 = 1021945
$ws.Range("B400").Value2 = 'Sat, Mar 22, 2025'
$ws.Range("C400").Value2 = 'Tri-City Americans'
$ws.Range("D400").Value2 = 'Spokane Chiefs'
$ws.Range("E400").Value2 = 'Spokane Chiefs'
$ws.Range("F400").Value2 = 'Spokane Chiefs'
$ws.Range("G400").Value2 = 1

$ws.Range("A401").Value2 = 1021948
$ws.Range("B401").Value2 = 'Sun, Mar 23, 2025'
$ws.Range("C401").Value2 = 'Edmonton Oil Kings'
$ws.Range("D401").Value2 = 'Red Deer Rebels'
$ws.Range("E401").Value2 = 'Edmonton Oil Kings'
$ws.Range("F401").Value2 = 'Edmonton Oil Kings'
$ws.Range("G401").Value2 = 1

$ws.Range("A402").Value2 = 1021947
$ws.Range("B402").Value2 = 'Sun, Mar 23, 2025'
$ws.Range("C402").Value2 = 'Calgary Hitmen'
$ws.Range("D402").Value2 = 'Medicine Hat Tigers'
$ws.Range("E402").Value2 = 'Medicine Hat Tigers'
$ws.Range("F402").Value2 = 'Medicine Hat Tigers'
$ws.Range("G402").Value2 = 1

$ws.Range("A403").Value2 = 1021949
$ws.Range("B403").Value2 = 'Sun, Mar 23, 2025'
$ws.Range("C403").Value2 = 'Vancouver Giants'
$ws.Range("D403").Value2 = 'Kelowna Rockets'
$ws.Range("E403").Value2 = 'Vancouver Giants'
$ws.Range("F403").Value2 = 'Vancouver Giants'
$ws.Range("G403").Value2 = 1

$ws.Range("A404").Value2 = 1021891
$ws.Range("B404").Value2 = 'Sun, Mar 23, 2025'
$ws.Range("C404").Value2 = 'Tri-City Americans'
$ws.Range("D404").Value2 = 'Seattle Thunderbirds'
$ws.Range("E404").Value2 = 'Seattle Thunderbirds'
$ws.Range("F404").Value2 = 'Seattle Thunderbirds'
$ws.Range("G404").Value2 = 1

$ws.Range("A405").Value2 = 1021963
$ws.Range("B405").Value2 = 'Fri, Mar 28, 2025'
$ws.Range("C405").Value2 = 'Calgary Hitmen'
$ws.Range("D405").Value2 = 'Saskatoon Blades'
$ws.Range("E405").Value2 = 'Saskatoon Blades'

$ws.Range("A406").Value2 = 1021965
$ws.Range("B406").Value2 = 'Fri, Mar 28, 2025'
$ws.Range("C406").Value2 = 'Lethbridge Hurricanes'
$ws.Range("D406").Value2 = 'Brandon Wheat Kings'
$ws.Range("E406").Value2 = 'Brandon Wheat Kings'

$ws.Range("A407").Value2 = 1021990
$ws.Range("B407").Value2 = 'Fri, Mar 28, 2025'
$ws.Range("C407").Value2 = 'Prince Albert Raiders'
$ws.Range("D407").Value2 = 'Edmonton Oil Kings'
$ws.Range("E407").Value2 = 'Prince Albert Raiders'

$ws.Range("A408").Value2 = 1022005
$ws.Range("B408").Value2 = 'Fri, Mar 28, 2025'
$ws.Range("C408").Value2 = 'Medicine Hat Tigers'
$ws.Range("D408").Value2 = 'Swift Current Broncos'
$ws.Range("E408").Value2 = 'Medicine Hat Tigers'

$ws.Range("A409").Value2 = 1021976
$ws.Range("B409").Value2 = 'Fri, Mar 28, 2025'
$ws.Range("C409").Value2 = 'Prince George Cougars'
$ws.Range("D409").Value2 = 'Portland Winterhawks'
$ws.Range("E409").Value2 = 'Portland Winterhawks'

$ws.Range("A410").Value2 = 1021997
$ws.Range("B410").Value2 = 'Fri, Mar 28, 2025'
$ws.Range("C410").Value2 = 'Vancouver Giants'
$ws.Range("D410").Value2 = 'Spokane Chiefs'
$ws.Range("E410").Value2 = 'Spokane Chiefs'

$ws.Range("A411").Value2 = 1021961
$ws.Range("B411").Value2 = 'Fri, Mar 28, 2025'
$ws.Range("C411").Value2 = 'Everett Silvertips'
$ws.Range("D411").Value2 = 'Seattle Thunderbirds'
$ws.Range("E411").Value2 = 'Seattle Thunderbirds'

$ws.Range("A412").Value2 = 1021974
$ws.Range("B412").Value2 = 'Fri, Mar 28, 2025'
$ws.Range("C412").Value2 = 'Victoria Royals'
$ws.Range("D412").Value2 = 'Tri-City Americans'
$ws.Range("E412").Value2 = 'Victoria Royals'

# Refresh the view: scroll to the bottom of the newly extended table and
# move the active selection to the next empty prediction row.
$ws.Range("A388").Select() | Out-Null
$ws.Range("D418").Select() | Out-Null
